$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("M2").Value = 2.815739333333333
$ws.Range("N2").Value = 8.447217999999999
$ws.Range("O2").Value = 0.07700398964630729
$ws.Range("P2").Value = 0.07700398964630729
$ws.Range("Q2").Value = 0.5927844617297777
$ws.Range("R2").Value = 5.335060155568
$ws.Range("S2").Value = 0.07700398964630729
$ws.Range("T2").Value = 0.07700398964630729

# Row 3
$ws.Range("O3").Value = 0.1324338085883186
$ws.Range("P3").Value = 0.1324338085883186
$ws.Range("S3").Value = 0.1324338085883186
$ws.Range("T3").Value = 0.1324338085883186

# Row 4
$ws.Range("M4").Value = 5.537790999999999
$ws.Range("N4").Value = 16.613373
$ws.Range("O4").Value = 0.1514458372546134
$ws.Range("P4").Value = 0.1514458372546134
$ws.Range("Q4").Value = 1.165845296205333
$ws.Range("R4").Value = 10.492607665848
$ws.Range("S4").Value = 0.1514458372546134
$ws.Range("T4").Value = 0.1514458372546134

# Row 5
$ws.Range("M5").Value = 1.188595666666667
$ws.Range("N5").Value = 3.565787
$ws.Range("O5").Value = 0.03250535563648733
$ws.Range("P5").Value = 0.03250535563648733
$ws.Range("Q5").Value = 0.2502294989235556
$ws.Range("R5").Value = 2.252065490312
$ws.Range("S5").Value = 0.03250535563648733
$ws.Range("T5").Value = 0.03250535563648733

# Row 6
$ws.Range("M6").Value = 18.85109966666667
$ws.Range("N6").Value = 56.553299
$ws.Range("O6").Value = 0.5155341854158992
$ws.Range("P6").Value = 0.5155341854158992
$ws.Range("Q6").Value = 3.968634041024889
$ws.Range("R6").Value = 35.717706369224
$ws.Range("S6").Value = 0.5155341854158992
$ws.Range("T6").Value = 0.5155341854158992

# Row 7
$ws.Range("M7").Value = 3.330328666666666
$ws.Range("N7").Value = 9.990985999999999
$ws.Range("O7").Value = 0.09107682345837424
$ws.Range("P7").Value = 0.09107682345837424
$ws.Range("Q7").Value = 0.7011185526595556
$ws.Range("R7").Value = 6.310066973936
$ws.Range("S7").Value = 0.09107682345837424
$ws.Range("T7").Value = 0.09107682345837424
